$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.262.63"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").Value = "1.857.13"
$ws.Range("E3").Value = "  +0.30%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'0.7008"
$ws.Range("E5").Value = "  +2.13%  "
$ws.Range("D6").Value = "'238.06"
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "'0.08028"
$ws.Range("E8").Value = "  +3.17%  "
$ws.Range("D9").Value = "'0.3021"
$ws.Range("E9").Value = "  -0.51%  "
$ws.Range("D10").Value = "'23.49"
$ws.Range("E10").Value = "  +1.29%  "
$ws.Range("D11").Value = "'0.08189"
$ws.Range("E11").Value = "  +0.32%  "
$ws.Range("D12").Value = "1.868.27"
$ws.Range("E12").Value = "  +1.40%  "
$ws.Range("D13").Value = "'5.206"
$ws.Range("E13").Value = "  +0.12%  "
$ws.Range("D14").Value = "'0.7070"
$ws.Range("E14").Value = "  -2.29%  "
$ws.Range("D15").Value = "'89.68"
$ws.Range("E15").Value = "  +0.58%  "
$ws.Range("D16").Value = "29.309.54"
$ws.Range("E16").Value = "  +0.57%  "
$ws.Range("D17").Value = "'5.827"
$ws.Range("E17").Value = "  +1.69%  "
$ws.Range("D18").Value = "'0.000007895"
$ws.Range("E18").Value = "  +1.06%  "
$ws.Range("D19").Value = "'13.28"
$ws.Range("E19").Value = "  +0.94%  "
$ws.Range("D20").Value = "'237.64"
$ws.Range("E20").Value = "  +1.62%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "'1.001"
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "2.126.71"
$ws.Range("E22").Value = "  +0.98%  "
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("E24").Value = "  -0.83%  "
$ws.Range("D25").Value = "'163.11"
$ws.Range("E25").Value = "  +0.88%  "
$ws.Range("D26").Value = "'8.893"
$ws.Range("E26").Value = "  -0.72%  "
$ws.Range("D27").Value = "'0.1425"
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("D28").Value = "'18.10"
$ws.Range("E28").Value = "  +0.18%  "
$ws.Range("D29").Value = "'1.920"
$ws.Range("E29").Value = "  -2.25%  "
$ws.Range("D30").Value = "'1.412"
$ws.Range("E30").Value = "  +0.51%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'1.477"
$ws.Range("E31").Value = "  -0.37%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'4.374"
$ws.Range("E32").Value = "  -3.19%  "
$ws.Range("D33").Value = "'4.038"
$ws.Range("E33").Value = "  +0.83%  "
$ws.Range("D34").Value = "'0.05194"
$ws.Range("E34").Value = "  +0.14%  "
$ws.Range("D35").Value = "'1.164"
$ws.Range("E35").Value = "  -1.03%  "
$ws.Range("D36").Value = "'0.7210"
$ws.Range("E36").Value = "  +2.52%  "
$ws.Range("D37").Value = "'1.001"
$ws.Range("E37").Value = "  -2.87%  "
$ws.Range("D38").Value = "'2.703"
$ws.Range("E38").Value = "  +1.84%  "
$ws.Range("D39").Value = "'0.01850"
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("E40").Value = "  +1.73%  "
$ws.Range("D41").Value = "'0.9383"
$ws.Range("E41").Value = "  +2.36%  "
$ws.Range("D42").Value = "1.152.59"
$ws.Range("E42").Value = "  +4.81%  "
$ws.Range("D43").Value = "'5.997"
$ws.Range("E43").Value = "  -0.26%  "
$ws.Range("D44").Value = "'0.4264"
$ws.Range("E44").Value = "  -0.33%  "
$ws.Range("D45").Value = "'70.46"
$ws.Range("E45").Value = "  +0.28%  "
$ws.Range("E46").Value = "  +0.09%  "
$ws.Range("D47").Value = "'102.89"
$ws.Range("E47").Value = "  +0.50%  "
$ws.Range("D48").Value = "'0.5286"
$ws.Range("E48").Value = "  -3.91%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "'1.759"
$ws.Range("E49").Value = "  +0.29%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.023.22"
$ws.Range("E50").Value = "  +1.37%  "
$ws.Range("D51").Value = "'9.159"
$ws.Range("E51").Value = "  +0.22%  "
